$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-375 all hold the serial date 45175
# (2023-09-06) and must be updated to 45177 (2023-09-08).
$range = $ws.Range("C2:C375")
$range.Value = 45177
